# Applies the commit "Mudei a ordem do Local mais perto da escolha do
# tupo de TCC." which:
#   1. Swaps the "Local:" label (first row) with "Eixo:" (becomes "Eixo:")
#   2. Swaps the "Eixo:" label (second row) with "Local:" (becomes "Local:")
#   3. Adds a "< " before "2.000 caracteres)" in the Resumo do Problema item

$d = $word.ActiveDocument

# Swap the two labels without one replacement colliding with the other:
# stage the first through a placeholder token that cannot already occur
# in the document.
$placeholder = "@@SWAP_TOKEN@@"

# 1a) "Local:" -> placeholder
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Local:", $true, $false, $false, $false, $false, `
    $true, 1, $false, $placeholder, 2)

# 2) "Eixo:" -> "Local:"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Eixo:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Local:", 2)

# 1b) placeholder -> "Eixo:"
$rng1b = $d.Content
$found1b = $rng1b.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
    $true, 1, $false, "Eixo:", 2)

# 3) "o Problema (2.000 caracteres)" -> "o Problema (< 2.000 caracteres)"
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("o Problema (2.000 caracteres)", $true, $false, $false, `
    $false, $false, $true, 1, $false, "o Problema (< 2.000 caracteres)", 2)

Write-Host "Local->token:" $found1 " Eixo->Local:" $found2 " token->Eixo:" $found1b " Resumo:" $found3
